# Apply the "feat: add 2022-Q3 data" edit to the workbook.
#
# Summary of the change:
#  1. A brand-new sheet "2022-Q3" is inserted right after "总计" and before
#     "2022-Q2". It carries the same layout/format as the other quarter
#     sheets (it is a duplicate of "2022-Q2", whose values are then
#     overwritten with the new Q3 numbers).
#  2. The "总计" (summary) sheet gets a new first data row for 2022-Q3 and
#     every existing row shifts down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating the "2022-Q2" sheet
#    (this keeps column layout, headers, number formats, etc. identical
#    to the other quarterly detail sheets) and placing the copy in front
#    of it.
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)

# The copy is inserted immediately before "2022-Q2" and is named
# "2022-Q2 (2)" by default - grab it and rename it.
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

# Overwrite the data row with the new 2022-Q3 figures. Columns D-G keep
# their original "text" formatting (the source data stores these ratios
# as text), column H is a plain number.
$wsQ3.Range("D2:G2").NumberFormat = "@"
$wsQ3.Range("D2").Value = "46.95"
$wsQ3.Range("E2").Value = "92.53"
$wsQ3.Range("F2").Value = "3.85"
$wsQ3.Range("G2").Value = "1.8076"
$wsQ3.Range("H2").Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" overview sheet: insert the 2022-Q3 row at the top
#    of the data and push the rest of the rows down by one.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Copy the formatting of an existing data cell in column A down onto the
# newly used row 8 before the values are written, so it matches the
# styling (border/alignment) of the other index cells in column A.
$wsTotal.Range("A7").Copy()
$wsTotal.Range("A8").PasteSpecial(-4122)

$wsTotal.Range("A8").Value = 6
$wsTotal.Range("B8").Value = "2020-Q4"
$wsTotal.Range("C8").Value = 6
$wsTotal.Range("D8").Value = 1.15

$wsTotal.Range("A7").Value = 5
$wsTotal.Range("B7").Value = "2021-Q1"
$wsTotal.Range("C7").Value = 2
$wsTotal.Range("D7").Value = 0.24

$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = "2021-Q2"
$wsTotal.Range("C6").Value = 3
$wsTotal.Range("D6").Value = 0.3

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2021-Q3"
$wsTotal.Range("C5").Value = 3
$wsTotal.Range("D5").Value = 2.19

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q4"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 2.52

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 1.58

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 1.81
